# Employee processing pipeline update:
#  - phone_number (E2) switches from a text value to a real number
#  - nine additional employee rows (3-12) are appended below the existing one
#  - the very last row (12) keeps its phone number as text (data-quality quirk
#    carried over from the source pipeline), matching the rest of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: convert the existing phone number cell from text to a number ---
$ws.Range("E2").Value = 6145551234

# --- New rows 3-12 ---
$rows = @(
    @{ First = "Jessica"; Last = "Martinez"; Dob = "1992-11-05"; Email = "j.martinez@email.com"; Phone = 6195559876; PhoneIsText = $false; Position = "marketing specialist" },
    @{ First = "Jessica"; Last = "Martinez"; Dob = "1992-11-05"; Email = "j.martinez@email.com"; Phone = 6195559876; PhoneIsText = $false; Position = "marketing specialist" },
    @{ First = "Natalie"; Last = "Brown";    Dob = "1988-05-29"; Email = "n.brown@email.com";    Phone = 3035552468; PhoneIsText = $false; Position = "administrative assistant" },
    @{ First = "Michael"; Last = "Johnson";  Dob = "1990-07-22"; Email = "m.johnson@email.com";  Phone = 6145551234; PhoneIsText = $false; Position = "customer service rep" },
    @{ First = "Jessica"; Last = "Martinez"; Dob = "1992-11-05"; Email = "j.martinez@email.com"; Phone = 6195559876; PhoneIsText = $false; Position = "marketing specialist" },
    @{ First = "Natalie"; Last = "Brown";    Dob = "1988-05-29"; Email = "n.brown@email.com";    Phone = 3035552468; PhoneIsText = $false; Position = "administrative assistant" },
    @{ First = "Robert";  Last = "Anderson"; Dob = "1976-09-10"; Email = "r.anderson@email.com"; Phone = 6125556789; PhoneIsText = $false; Position = "project manager" },
    @{ First = "Robert";  Last = "Anderson"; Dob = "1976-09-10"; Email = "r.anderson@email.com"; Phone = 6125556789; PhoneIsText = $false; Position = "project manager" },
    @{ First = "Robert";  Last = "Anderson"; Dob = "1976-09-10"; Email = "r.anderson@email.com"; Phone = 6125556789; PhoneIsText = $false; Position = "project manager" },
    @{ First = "Robert";  Last = "Anderson"; Dob = "1976-09-10"; Email = "r.anderson@email.com"; Phone = "6125556789"; PhoneIsText = $true;  Position = "project manager" }
)

$r = 3
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.First
    $ws.Range("B$r").Value = $row.Last

    # Text that looks like a date must be forced to stay text, otherwise Excel
    # auto-converts it into a date serial number.
    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $row.Dob
    $ws.Range("C$r").Style = "Normal"

    $ws.Range("D$r").Value = $row.Email

    if ($row.PhoneIsText) {
        # Keep this one as a genuine text value (digits-only strings are
        # otherwise auto-coerced into numbers by Excel on assignment).
        $ws.Range("E$r").NumberFormat = "@"
        $ws.Range("E$r").Value = $row.Phone
        $ws.Range("E$r").Style = "Normal"
    } else {
        $ws.Range("E$r").Value = $row.Phone
    }

    $ws.Range("F$r").Value = $row.Position

    $r = $r + 1
}
